$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 104. This pushes the existing rows
# 104-182 down to 106-184, matching the target layout.
$ws.Rows("104:105").Insert()

# Populate the two newly inserted rows with the new records.
# Row 104
$ws.Cells.Item(104, 1).Value = 5
$ws.Cells.Item(104, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(104, 3).Value = "Maule"
$ws.Cells.Item(104, 4).Value = 45126
$ws.Cells.Item(104, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(104, 5).Value = 7
$ws.Cells.Item(104, 6).Value = 100112001
$ws.Cells.Item(104, 7).Value = "Berenjena"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 200
$ws.Cells.Item(104, 11).Value = 7000
$ws.Cells.Item(104, 12).Value = 7000
$ws.Cells.Item(104, 13).Value = 7000
$ws.Cells.Item(104, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(104, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(104, 16).Value = 140
$ws.Cells.Item(104, 17).Value = 50
$ws.Cells.Item(104, 18).Value = "Hortaliza"

# Row 105
$ws.Cells.Item(105, 1).Value = 5
$ws.Cells.Item(105, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(105, 3).Value = "Maule"
$ws.Cells.Item(105, 4).Value = 45126
$ws.Cells.Item(105, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(105, 5).Value = 7
$ws.Cells.Item(105, 6).Value = 100112001
$ws.Cells.Item(105, 7).Value = "Berenjena"
$ws.Cells.Item(105, 8).Value = "Sin especificar"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 300
$ws.Cells.Item(105, 11).Value = 7000
$ws.Cells.Item(105, 12).Value = 7000
$ws.Cells.Item(105, 13).Value = 7000
$ws.Cells.Item(105, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(105, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(105, 16).Value = 140
$ws.Cells.Item(105, 17).Value = 50
$ws.Cells.Item(105, 18).Value = "Hortaliza"
